$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bio-parameters")

# Update Thunder Bay (Lake Superior) spawn temperature data
$ws.Range("I11").Value = 4.943492
$ws.Range("J11").Value = 4.3679399999999999

# Update the selected cell shown in the saved workbook
$ws.Range("I18").Select()
